$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("G:H").Insert()

$ws.Range("G1").Value = "Logo URL"
$ws.Range("H1").Value = "Website"

$ws.Range("G2").Value = "https://www.skiresort.info/fileadmin/_processed_/84/42/2d/d8/0326e4726a.png"
$ws.Range("H2").Value = "https://www.kitzski.at/"

$ws.Range("G3").Value = "https://www.skiresort.info/fileadmin/_processed_/b7/77/72/27/3b74a66c33.png"
$ws.Range("H3").Value = "https://www.zillertalarena.com"
